# Mobile App translation logic
# Fill in the "Assigned to" column (J) for several rows, populate a new
# table row (21) with Vehicle / Service area Add/update/delete / P1, and
# move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry for the previously-blank row 21
$ws.Range("B21").Value = "Vehicle"
$ws.Range("C21").Value = "Service area Add/update/delete"
$ws.Range("F21").Value = "P1"

# New "assigned to" values in column J for various rows
$ws.Range("J4").Value = "Deepak"
$ws.Range("J5").Value = "Ravi"
$ws.Range("J9").Value = "Atul"
$ws.Range("J12").Value = "Viki"
$ws.Range("J10").Value = "MP"

# Move the selection as it ended up after the edits
$ws.Range("L23").Select() | Out-Null
